$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Nathanael%Lapidus%nathanael.lapidus@inserm.fr%1,                    Xianlong%Zhou%xianlongzhou@whu.edu.cn%2,                    Xianlong%Zhou%xianlongzhou@whu.edu.cn%0,                    Fabrice%Carrat%fabrice.carrat@iplesp.upmc.fr%2,                    Fabrice%Carrat%fabrice.carrat@iplesp.upmc.fr%0,                    Bruno%Riou%bruno.riou@aphp.fr%2,                    Bruno%Riou%bruno.riou@aphp.fr%0,                    Yan%Zhao%doctoryanzhao@whu.edu.cn%0,                    Yan%Zhao%doctoryanzhao@whu.edu.cn%0,                    Gilles%Hejblum%gilles.hejblum@inserm.fr%2,                    Gilles%Hejblum%gilles.hejblum@inserm.fr%0]'
$ws.Range("I2").Value = '_PMC_Springer_CROSSREF'
$ws.Range("E3").Value = '[Xiaobo%Feng%NULL%0,                    Peiyun%Li%NULL%1,                    Liang%Ma%NULL%1,                    Hang%Liang%NULL%1,                    Jie%Lei%NULL%0,                    Wenqiang%Li%NULL%1,                    Kun%Wang%NULL%0,                    Yu%Song%NULL%1,                    Shuai%Li%NULL%1,                    Wei%Yang%NULL%1,                    Cao%Yang%NULL%1]'
$ws.Range("I3").Value = '_PMC_CROSSREF'
$ws.Range("E4").Value = '[ Q.%Bi%null%1,                     C.% Hong%null%1,                     J.% Meng%null%1,                     Z.% Wu%null%1,                     P.% Zhou%null%1,                     C.% Ye%null%1,                     B.% Sun%null%1,                     L. M.% Kucirka%null%1,                     A. S.% Azman%null%1,                     T.% Wang%null%0,                     J.% Chen%null%1,                     Z.% Wang%null%1,                     L.% Liu%null%0,                     J.% Lessler%null%1,                     J. K.% Edwards%null%1,                     T.% Ma%null%1,                     G. % Zhang%null%1]'
$ws.Range("I4").Value = '_MedBiorxiv_CROSSREF'
$ws.Range("E5").Value = '[Fei%Zhou%NULL%0,                    Ting%Yu%NULL%0,                    Ronghui%Du%NULL%0,                    Guohui%Fan%NULL%0,                    Ying%Liu%NULL%0,                    Zhibo%Liu%NULL%0,                    Jie%Xiang%NULL%0,                    Yeming%Wang%NULL%0,                    Bin%Song%NULL%0,                    Xiaoying%Gu%NULL%0,                    Lulu%Guan%NULL%0,                    Yuan%Wei%NULL%0,                    Hui%Li%NULL%0,                    Xudong%Wu%NULL%0,                    Jiuyang%Xu%NULL%0,                    Shengjin%Tu%NULL%0,                    Yi%Zhang%NULL%0,                    Hua%Chen%NULL%0,                    Bin%Cao%NULL%0]'
$ws.Range("I5").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E6").Value = '[Yonghao%Xu%NULL%1,                    Zhiheng%Xu%NULL%1,                    Xuesong%Liu%NULL%1,                    Lihua%Cai%NULL%1,                    Haichong%Zheng%NULL%1,                    Yongbo%Huang%NULL%1,                    Lixin%Zhou%NULL%1,                    Linxi%Huang%NULL%1,                    Yun%Ling%NULL%0,                    Liehua%Deng%NULL%1,                    Jianwei%Li%NULL%0,                    Sibei%Chen%NULL%1,                    Dongdong%Liu%NULL%1,                    Zhimin%Lin%NULL%1,                    Liang%Zhou%NULL%1,                    Weiqun%He%NULL%1,                    Nanshan%Zhong%NULL%0,                    Xiaoqing%Liu%NULL%3,                    Yimin%Li%NULL%0]'
$ws.Range("I6").Value = '_PMC_CROSSREF'
$ws.Range("E7").Value = '[Yi%Zheng%NULL%0,                    Li-jun%Sun%NULL%1,                    Mi%Xu%NULL%1,                    Jian%Pan%NULL%1,                    Yun-tao%Zhang%NULL%1,                    Xue-ling%Fang%NULL%1,                    Qiang%Fang%NULL%1,                    Hong-liu%Cai%NULL%1]'
$ws.Range("I7").Value = '_PMC_Springer_CROSSREF'
$ws.Range("E8").Value = '[Rong-Hui%Du%NULL%0,                    Li-Min%Liu%NULL%2,                    Wen%Yin%NULL%0,                    Wen%Wang%NULL%0,                    Lu-Lu%Guan%NULL%2,                    Ming-Li%Yuan%NULL%0,                    Yu-Lei%Li%NULL%2,                    Yi%Hu%NULL%0,                    Xu-Yan%Li%NULL%0,                    Bing%Sun%NULL%2,                    Peng%Peng%NULL%0,                    Huan-Zhong%Shi%NULL%0]'
$ws.Range("I8").Value = '_PMC_CROSSREF'
$ws.Range("E9").Value = '[Jianlei%Cao%NULL%0,                    Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%0,                    Wenlin%Cheng%NULL%0,                    Lei%Yu%NULL%0,                    Ya-Kun%Liu%NULL%0,                    Xiaoyong%Hu%NULL%0,                    Qiang%Liu%NULL%0]'
$ws.Range("I9").Value = '_PMC_CROSSREF'
$ws.Range("E10").Value = '[Wei-jie%Guan%NULL%0,                    Zheng-yi%Ni%NULL%0,                    Zheng-yi%Ni%NULL%0,                    Yu%Hu%NULL%0,                    Wen-hua%Liang%NULL%0,                    Chun-quan%Ou%NULL%0,                    Jian-xing%He%NULL%0,                    Lei%Liu%NULL%0,                    Hong%Shan%NULL%0,                    Chun-liang%Lei%NULL%0,                    David S.C.%Hui%NULL%0,                    Bin%Du%NULL%0,                    Lan-juan%Li%NULL%0,                    Guang%Zeng%NULL%0,                    Kwok-Yung%Yuen%NULL%0,                    Ru-chong%Chen%NULL%0,                    Chun-li%Tang%NULL%0,                    Tao%Wang%NULL%0,                    Ping-yan%Chen%NULL%0,                    Jie%Xiang%NULL%0,                    Shi-yue%Li%NULL%0,                    Jin-lin%Wang%NULL%0,                    Zi-jing%Liang%NULL%0,                    Yi-xiang%Peng%NULL%0,                    Li%Wei%NULL%0,                    Yong%Liu%NULL%0,                    Ya-hua%Hu%NULL%0,                    Peng%Peng%NULL%0,                    Jian-ming%Wang%NULL%0,                    Ji-yang%Liu%NULL%0,                    Zhong%Chen%NULL%0,                    Gang%Li%NULL%0,                    Zhi-jian%Zheng%NULL%0,                    Shao-qin%Qiu%NULL%0,                    Jie%Luo%NULL%0,                    Chang-jiang%Ye%NULL%0,                    Shao-yong%Zhu%NULL%0,                    Nan-shan%Zhong%NULL%0]'
$ws.Range("I10").Value = '_PMC_CROSSREF'
$ws.Range("E11").Value = '[Shaoqing%Lei%NULL%0,                    Fang%Jiang%NULL%0,                    Wating%Su%NULL%2,                    Chang%Chen%NULL%2,                    Jingli%Chen%NULL%2,                    Wei%Mei%NULL%2,                    Li-Ying%Zhan%NULL%2,                    Yifan%Jia%NULL%2,                    Liangqing%Zhang%NULL%2,                    Danyong%Liu%NULL%2,                    Zhong-Yuan%Xia%xiazhongyuan2005@aliyun.com%2,                    Zhengyuan%Xia%zyxia@hkucc.hku.hk%0]'
$ws.Range("I11").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E12").Value = '[Xiaobo%Yang%NULL%0,                    Yuan%Yu%NULL%0,                    Jiqian%Xu%NULL%0,                    Huaqing%Shu%NULL%0,                    Jia''an%Xia%NULL%0,                    Hong%Liu%NULL%0,                    Yongran%Wu%NULL%0,                    Lu%Zhang%NULL%0,                    Zhui%Yu%NULL%0,                    Minghao%Fang%NULL%0,                    Ting%Yu%NULL%0,                    Yaxin%Wang%NULL%0,                    Shangwen%Pan%NULL%0,                    Xiaojing%Zou%NULL%0,                    Shiying%Yuan%NULL%0,                    You%Shang%NULL%0]'
$ws.Range("I12").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E13").Value = '[Chaolin%Huang%NULL%0,                    Yeming%Wang%NULL%0,                    Xingwang%Li%NULL%0,                    Lili%Ren%NULL%0,                    Jianping%Zhao%NULL%0,                    Yi%Hu%NULL%0,                    Li%Zhang%NULL%0,                    Guohui%Fan%NULL%0,                    Jiuyang%Xu%NULL%0,                    Xiaoying%Gu%NULL%0,                    Zhenshun%Cheng%NULL%0,                    Ting%Yu%NULL%0,                    Jiaan%Xia%NULL%0,                    Yuan%Wei%NULL%0,                    Wenjuan%Wu%NULL%0,                    Xuelei%Xie%NULL%0,                    Wen%Yin%NULL%0,                    Hui%Li%NULL%0,                    Min%Liu%NULL%0,                    Yan%Xiao%NULL%0,                    Hong%Gao%NULL%0,                    Li%Guo%NULL%0,                    Jungang%Xie%NULL%0,                    Guangfa%Wang%NULL%0,                    Rongmeng%Jiang%NULL%0,                    Zhancheng%Gao%NULL%0,                    Qi%Jin%NULL%0,                    Jianwei%Wang%wangjw28@163.com%0,                    Bin%Cao%caobin_ben@163.com%0]'
$ws.Range("I13").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E14").Value = '[Yang%Wang%NULL%1,                    Xiaofan%Lu%NULL%1,                    Yongsheng%Li%NULL%2,                    Yongsheng%Li%NULL%0,                    Hui%Chen%NULL%3,                    Taige%Chen%NULL%1,                    Nan%Su%NULL%1,                    Fang%Huang%NULL%1,                    Jing%Zhou%NULL%1,                    Bing%Zhang%NULL%1,                    Fangrong%Yan%NULL%2,                    Fangrong%Yan%NULL%0,                    Jun%Wang%NULL%2,                    Jun%Wang%NULL%0]'
$ws.Range("I14").Value = '_PMC_CROSSREF'
$ws.Range("E15").Value = '[Antonelli%Massimo%coreGivesNoEmail%0,                  Cabrini%Luca%coreGivesNoEmail%0,                  Castelli%Antonio%coreGivesNoEmail%0,                  Cecconi%Maurizio%coreGivesNoEmail%0,                  Cereda%Danilo%coreGivesNoEmail%0,                  Coluccello%Antonio%coreGivesNoEmail%0,                  Foti%Giuseppe%coreGivesNoEmail%0,                  Fumagalli%Roberto%coreGivesNoEmail%0,                  Grasselli%Giacomo%coreGivesNoEmail%0,                  Iotti%Giorgio%coreGivesNoEmail%0,                  Latronico%Nicola%coreGivesNoEmail%0,                  Lorini%Luca%coreGivesNoEmail%0,                  Merler%Stefano%coreGivesNoEmail%0,                  Natalini%Giuseppe%coreGivesNoEmail%0,                  Pesenti%Antonio%coreGivesNoEmail%0,                  Piatti%Alessandra%coreGivesNoEmail%0,                  Ranieri%Marco Vito%coreGivesNoEmail%0,                  Scandroglio%Anna Mara%coreGivesNoEmail%0,                  Storti%Enrico%coreGivesNoEmail%0,                  Zanella%Alberto%coreGivesNoEmail%0,                  Zangrillo%Alberto%coreGivesNoEmail%0]'
$ws.Range("I15").Value = '_CROSSREF'
$ws.Range("E16").Value = '[Arthur%Simonnet%NULL%0,                    Mikael%Chetboun%NULL%1,                    Julien%Poissy%NULL%1,                    Violeta%Raverdy%NULL%1,                    Jerome%Noulette%NULL%2,                    Jerome%Noulette%NULL%0,                    Alain%Duhamel%NULL%0,                    Julien%Labreuche%NULL%0,                    Daniel%Mathieu%NULL%2,                    Francois%Pattou%francois.pattou@univ-lille.fr%1,                    Merce%Jourdain%mercedes.jourdain@univ-lille.fr%2,                    Merce%Jourdain%mercedes.jourdain@univ-lille.fr%0,                    Robert%Caizzo%NULL%2,                    Robert%Caizzo%NULL%0,                    Morgan%Caplan%NULL%1,                    Nicolas%Cousin%NULL%1,                    Thibault%Duburcq%NULL%1,                    Arthur%Durand%NULL%1,                    Ahmed%El kalioubie%NULL%1,                    Raphael%Favory%NULL%1,                    Bruno%Garcia%NULL%1,                    Patrick%Girardie%NULL%1,                    Julien%Goutay%NULL%1,                    Marion%Houard%NULL%1,                    Emmanuelle%Jaillette%NULL%1,                    Nicolas%Kostuj%NULL%1,                    Geoffrey%Ledoux%NULL%1,                    Daniel%Mathieu%NULL%0,                    Anne Sophie%Moreau%NULL%1,                    Christopher%Niles%NULL%1,                    Saad%Nseir%NULL%1,                    Thierry%Onimus%NULL%1,                    Erika%Parmentier%NULL%1,                    Sebastien%Préau%NULL%1,                    Laurent%Robriquet%NULL%1,                    Anahita%Rouze%NULL%1,                    Sophie%Six%NULL%1,                    Hélène%Verkindt%NULL%1]'
$ws.Range("I16").Value = '_PMC_CROSSREF'
$ws.Range("E17").Value = '[Matt%Arentz%xref no email%0,            Eric%Yim%xref no email%0,            Lindy%Klaff%xref no email%0,            Sharukh%Lokhandwala%xref no email%0,            Francis X.%Riedo%xref no email%0,            Maria%Chong%xref no email%0,            Melissa%Lee%xref no email%0]'
$ws.Range("I17").Value = '_CROSSREF'
$ws.Range("E18").Value = '[Michael G%Argenziano%NULL%0,                    Samuel L%Bruce%NULL%1,                    Cody L%Slater%NULL%1,                    Jonathan R%Tiao%NULL%1,                    Matthew R%Baldwin%NULL%2,                    R Graham%Barr%NULL%1,                    Bernard P%Chang%NULL%1,                    Katherine H%Chau%NULL%1,                    Justin J%Choi%NULL%1,                    Nicholas%Gavin%NULL%1,                    Parag%Goyal%NULL%1,                    Angela M%Mills%NULL%1,                    Ashmi A%Patel%NULL%1,                    Marie-Laure S%Romney%NULL%1,                    Monika M%Safford%NULL%1,                    Neil W%Schluger%NULL%1,                    Soumitra%Sengupta%NULL%1,                    Magdalena E%Sobieszczyk%NULL%1,                    Jason E%Zucker%NULL%1,                    Paul A%Asadourian%NULL%1,                    Fletcher M%Bell%NULL%1,                    Rebekah%Boyd%NULL%1,                    Matthew F%Cohen%NULL%1,                    MacAlistair I%Colquhoun%NULL%1,                    Lucy A%Colville%NULL%1,                    Joseph H%de Jonge%NULL%1,                    Lyle B%Dershowitz%NULL%1,                    Shirin A%Dey%NULL%1,                    Katherine A%Eiseman%NULL%1,                    Zachary P%Girvin%NULL%1,                    Daniella T%Goni%NULL%1,                    Amro A%Harb%NULL%1,                    Nicholas%Herzik%NULL%1,                    Sarah%Householder%NULL%1,                    Lara E%Karaaslan%NULL%1,                    Heather%Lee%NULL%1,                    Evan%Lieberman%NULL%1,                    Andrew%Ling%NULL%1,                    Ree%Lu%NULL%1,                    Arthur Y%Shou%NULL%1,                    Alexander C%Sisti%NULL%1,                    Zachary E%Snow%NULL%1,                    Colin P%Sperring%NULL%1,                    Yuqing%Xiong%NULL%1,                    Henry W%Zhou%NULL%1,                    Karthik%Natarajan%NULL%1,                    George%Hripcsak%NULL%1,                    Ruijun%Chen%NULL%1]'
$ws.Range("I18").Value = '_PMC_CROSSREF'
$ws.Range("E19").Value = '[Sara C.%Auld%NULL%0,                    Mark%Caridi-Scheible%NULL%0,                    James M.%Blum%NULL%0,                    Chad%Robichaux%NULL%0,                    Colleen%Kraft%NULL%0,                    Jesse T.%Jacob%NULL%0,                    Craig S.%Jabaley%NULL%0,                    David%Carpenter%NULL%0,                    Roberta%Kaplow%NULL%0,                    Alfonso C.%Hernandez-Romieu%NULL%0,                    Max W.%Adelman%NULL%0,                    Greg S.%Martin%NULL%0,                    Craig M.%Coopersmith%NULL%0,                    David J.%Murphy%NULL%0,                    NULL%NULL%NULL%0,                    NULL%NULL%NULL%0]'
$ws.Range("I19").Value = '_PMC_CROSSREF'
$ws.Range("E20").Value = '[Pavan K.%Bhatraju%NULL%0,                    Bijan J.%Ghassemieh%NULL%0,                    Michelle%Nichols%NULL%0,                    Richard%Kim%NULL%0,                    Keith R.%Jerome%NULL%0,                    Arun K.%Nalla%NULL%0,                    Alexander L.%Greninger%NULL%0,                    Sudhakar%Pipavath%NULL%0,                    Mark M.%Wurfel%NULL%0,                    Laura%Evans%NULL%0,                    Patricia A.%Kritek%NULL%0,                    T. Eoin%West%NULL%0,                    Andrew%Luks%NULL%0,                    Anthony%Gerbino%NULL%0,                    Chris R.%Dale%NULL%0,                    Jason D.%Goldman%NULL%0,                    Shane%O’Mahony%NULL%0,                    Carmen%Mikacenic%NULL%0]'
$ws.Range("I20").Value = '_PMC_CROSSREF'
$ws.Range("E21").Value = '[Laura C.%Myers%xref no email%1,            Stephen M.%Parodi%xref no email%1,            Gabriel J.%Escobar%xref no email%1,            Vincent X.%Liu%xref no email%1]'
$ws.Range("I21").Value = '_CROSSREF'
$ws.Range("E22").Value = '[Barnaby%D.P.%coreGivesNoEmail%0,                  Becker%L.B.%coreGivesNoEmail%0,                  Chelico%J.D.%coreGivesNoEmail%0,                  Cohen%S.L.%coreGivesNoEmail%0,                  Cookingham%J.%coreGivesNoEmail%0,                  Coppa%K.%coreGivesNoEmail%0,                  Crawford%J.M.%coreGivesNoEmail%0,                  Davidson%K.W.%coreGivesNoEmail%0,                  Diefenbach%M.A.%coreGivesNoEmail%0,                  Dominello%A.J.%coreGivesNoEmail%0,                  Duer-Hefele%J.%coreGivesNoEmail%0,                  Falzon%L.%coreGivesNoEmail%0,                  Gitlin%J.%coreGivesNoEmail%0,                  Hajizadeh%N.%coreGivesNoEmail%0,                  Harvin%T.G.%coreGivesNoEmail%0,                  Hirsch%J.S.%coreGivesNoEmail%0,                  Hirschwerk%D.A.%coreGivesNoEmail%0,                  Kim%E.J.%coreGivesNoEmail%0,                  Kozel%Z.M.%coreGivesNoEmail%0,                  Marrast%L.M.%coreGivesNoEmail%0,                  McGinn%T.%coreGivesNoEmail%0,                  Mogavero%J.N.%coreGivesNoEmail%0,                  Narasimhan%M.%coreGivesNoEmail%0,                  Osorio%G.A.%coreGivesNoEmail%0,                  Qiu%M.%coreGivesNoEmail%0,                  Richardson%S.%coreGivesNoEmail%0,                  Zanos%T.P.%coreGivesNoEmail%0]'
$ws.Range("I22").Value = '_CROSSREF'
$ws.Range("E23").Value = '[David R.%Ziehr%NULL%1,                    Jehan%Alladina%NULL%2,                    Jehan%Alladina%NULL%0,                    Camille R.%Petri%NULL%2,                    Camille R.%Petri%NULL%0,                    Jason H.%Maley%NULL%2,                    Jason H.%Maley%NULL%0,                    Ari%Moskowitz%NULL%1,                    Benjamin D.%Medoff%NULL%1,                    Kathryn A.%Hibbert%NULL%2,                    Kathryn A.%Hibbert%NULL%0,                    B. Taylor%Thompson%NULL%1,                    C. Corey%Hardin%NULL%1]'
$ws.Range("I23").Value = '_PMC_CROSSREF'
$ws.Range("E24").Value = '[Ishan%Paranjpe%NULL%0,                    Adam J%Russak%NULL%2,                    Adam J%Russak%NULL%0,                    Jessica K%De Freitas%NULL%1,                    Anuradha%Lala%NULL%1,                    Riccardo%Miotto%NULL%1,                    Akhil%Vaid%NULL%1,                    Kipp W%Johnson%NULL%1,                    Matteo%Danieletto%NULL%1,                    Eddye%Golden%NULL%1,                    Dara%Meyer%NULL%1,                    Manbir%Singh%NULL%1,                    Sulaiman%Somani%NULL%1,                    Arjun%Kapoor%NULL%1,                    Ross%O''Hagan%NULL%1,                    Sayan%Manna%NULL%1,                    Udit%Nangia%NULL%1,                    Suraj K%Jaladanki%NULL%1,                    Paul%O’Reilly%NULL%2,                    Paul%O’Reilly%NULL%0,                    Laura M%Huckins%NULL%1,                    Patricia%Glowe%NULL%1,                    Arash%Kia%NULL%1,                    Prem%Timsina%NULL%1,                    Robert M%Freeman%NULL%1,                    Matthew A%Levin%NULL%1,                    Jeffrey%Jhang%NULL%1,                    Adolfo%Firpo%NULL%1,                    Patricia%Kovatch%NULL%1,                    Joseph%Finkelstein%NULL%1,                    Judith A%Aberg%NULL%1,                    Emilia%Bagiella%NULL%1,                    Carol R%Horowitz%NULL%1,                    Barbara%Murphy%NULL%1,                    Zahi A%Fayad%NULL%1,                    Jagat%Narula%NULL%1,                    Eric J%Nestler%NULL%1,                    V%Fuster%NULL%1,                    Carlos%Cordon-Cardo%NULL%1,                    Dennis%Charney%NULL%1,                    David L%Reich%NULL%1,                    Allan%Just%NULL%1,                    Erwin P%Bottinger%NULL%1,                    Alexander W%Charney%NULL%1,                    Benjamin S%Glicksberg%NULL%1,                    Girish N%Nadkarni%NULL%1,                    NULL%NULL%NULL%0]'
$ws.Range("I24").Value = '_PMC_CROSSREF'
